$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows between each "Topic N" / keyword pair were previously empty gaps.
# Fill them in with descriptive summary text for each topic (no row shifting
# needed - row numbers of existing content stay the same).

$ws.Range("A3").Value = "Green and renewable power plant development instead of nuclear projets"
$ws.Range("A6").Value = "fueling electric and zero emission futuristic technology"
$ws.Range("A9").Value = "international event/summit/confernce on renewable energy and power sources"
$ws.Range("A12").Value = "Alternative fuel for NASA rocket launches and flights which creates carcinogenics"
$ws.Range("A15").Value = "State and government support towards climate change and clean renewable energy"
$ws.Range("A18").Value = "Clean byproducts and low carbon outcome"

$ws.Range("A13").Select()
